# Apply the review-flag correction to rows 22 and 23 (column G, "blue")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G22").Value = "no"
$ws.Range("G23").Value = "no"

# Reflect the saved selection state (active cell moved to A24)
$ws.Range("A24").Select()
